$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 94 (shifts existing rows 94..169 down to 95..170)
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across this dataset, so copy them
# from row 95 (the row that used to be row 94 before the insert).
$ws.Cells.Item(94, 1).Value2  = $ws.Cells.Item(95, 1).Value2
$ws.Cells.Item(94, 2).Value2  = $ws.Cells.Item(95, 2).Value2
$ws.Cells.Item(94, 3).Value2  = $ws.Cells.Item(95, 3).Value2
$ws.Cells.Item(94, 4).Value2  = 44589
$ws.Cells.Item(94, 5).Value2  = $ws.Cells.Item(95, 5).Value2
$ws.Cells.Item(94, 6).Value2  = $ws.Cells.Item(95, 6).Value2
$ws.Cells.Item(94, 7).Value2  = $ws.Cells.Item(95, 7).Value2
$ws.Cells.Item(94, 8).Value2  = $ws.Cells.Item(95, 8).Value2
$ws.Cells.Item(94, 9).Value2  = $ws.Cells.Item(95, 9).Value2
$ws.Cells.Item(94, 10).Value2 = 240
$ws.Cells.Item(94, 11).Value2 = 2500
$ws.Cells.Item(94, 12).Value2 = 3000
$ws.Cells.Item(94, 13).Value2 = 2750
$ws.Cells.Item(94, 14).Value2 = $ws.Cells.Item(95, 14).Value2
$ws.Cells.Item(94, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(94, 16).Value2 = 917
$ws.Cells.Item(94, 17).Value2 = $ws.Cells.Item(95, 17).Value2
$ws.Cells.Item(94, 18).Value2 = $ws.Cells.Item(95, 18).Value2
